# Add a new "verigen_add_verilog" worksheet, modeled on the existing
# "verigen_add_reference" sheet (same layout: Type/Description header,
# Prototype, Return value, Remarks, filename rows), positioned right
# after "verigen_add_reference" and before "read_excel_table".

$wb = $excel.ActiveWorkbook

$refSheet = $wb.Worksheets.Item("verigen_add_reference")

# Copy the reference sheet (preserves column widths, styles, row heights,
# merged layout, etc.) and place the copy right after it.
$refSheet.Copy($null, $refSheet)
$newSheet = $wb.Worksheets.Item($refSheet.Index + 1)
$newSheet.Name = "verigen_add_verilog"

# Update the prototype / remarks / parameter description text for the new
# function. Row 1 (Type/Description header), row 3 (Return value "-") stay
# as copied. (Written in this order so new shared-string entries land in
# the same index order as the authored workbook: Remarks text first, then
# the prototype, then the parameter description.)
$newSheet.Range("B4").Value = "Add verilog source(s)"
$newSheet.Range("B2").Value = "function verigen_add_verilog(filename)"
$newSheet.Range("B5").Value = "verilog file name (can use wildcard '*')"

# verigen_add_reference has an extra "desc" parameter row (row 6) that
# verigen_add_verilog does not need, since it only takes a filename.
$newSheet.Rows.Item(6).Delete()

# Match the printed page setup used by the other sheets.
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# The new sheet becomes the active / selected tab.
$newSheet.Range("A1").Select()
$newSheet.Activate()

# The old reference sheet is no longer the selected tab; its selection
# becomes the full used range instead of the single active cell it had.
$refSheet.Range("A1:B6").Select()
$refSheet.PageSetup.PaperSize = 9
$refSheet.PageSetup.Orientation = 1

# Re-activate the new sheet last so it ends up as the workbook's active tab.
$newSheet.Activate()
